# Hortaliza, Vega Modelo de Temuco - Acelga
# Two new daily price records are inserted right before the existing
# row that used to be row 475 (old dimension A1:R539). Every row from
# the old 475 onward shifts down by two rows, giving the new
# dimension A1:R541, and the two freshly inserted rows (now 475 and
# 476) carry the new "Fruta / hortaliza, semanal" data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 475-476; everything at/after the old row 475
# (old rows 475..539) moves down to 477..541.
$ws.Range("A475:A476").EntireRow.Insert()

# --- New row 475 ---
$ws.Range("A475").Value = 10
$ws.Range("B475").Value = "Vega Modelo de Temuco"
$ws.Range("C475").Value = "La Araucanía"
$ws.Range("D475").Value = 45127
$ws.Range("E475").Value = 9
$ws.Range("F475").Value = 100112009
$ws.Range("G475").Value = "Acelga"
$ws.Range("H475").Value = "Sin especificar"
$ws.Range("I475").Value = "Primera"
$ws.Range("J475").Value = 60
$ws.Range("K475").Value = 8000
$ws.Range("L475").Value = 8000
$ws.Range("M475").Value = 8000
$ws.Range("N475").Value = "$/docena de atados (12 kilos)"
$ws.Range("O475").Value = "Provincia de Cautín"
$ws.Range("P475").Value = 667
$ws.Range("Q475").Value = 12
$ws.Range("R475").Value = "Hortaliza"

# --- New row 476 ---
$ws.Range("A476").Value = 10
$ws.Range("B476").Value = "Vega Modelo de Temuco"
$ws.Range("C476").Value = "La Araucanía"
$ws.Range("D476").Value = 45127
$ws.Range("E476").Value = 9
$ws.Range("F476").Value = 100112009
$ws.Range("G476").Value = "Acelga"
$ws.Range("H476").Value = "Sin especificar"
$ws.Range("I476").Value = "Primera"
$ws.Range("J476").Value = 80
$ws.Range("K476").Value = 6000
$ws.Range("L476").Value = 6000
$ws.Range("M476").Value = 6000
$ws.Range("N476").Value = "$/docena de atados (12 kilos)"
$ws.Range("O476").Value = "Región Metropolitana"
$ws.Range("P476").Value = 500
$ws.Range("Q476").Value = 12
$ws.Range("R476").Value = "Hortaliza"
